$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record arrived for "Coco" (Mercado Mayorista Lo Valledor de
# Santiago). It belongs at the top of the date-ordered block (row 63), so insert
# a new row there and push the existing rows (63-88) down to (64-89).
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the new observation. The columns
# that stay constant for every row in this block (A, B, C, E, F, G, H, I, J, K, L, Q, T)
# are copied straight from the row below (which held this data before the insert).
$ws.Cells.Item(63, 1).Value = $ws.Cells.Item(64, 1).Value()
$ws.Cells.Item(63, 2).Value = $ws.Cells.Item(64, 2).Value()
$ws.Cells.Item(63, 3).Value = $ws.Cells.Item(64, 3).Value()
$ws.Cells.Item(63, 4).Value = 44845
$ws.Cells.Item(63, 5).Value = $ws.Cells.Item(64, 5).Value()
$ws.Cells.Item(63, 6).Value = $ws.Cells.Item(64, 6).Value()
$ws.Cells.Item(63, 7).Value = $ws.Cells.Item(64, 7).Value()
$ws.Cells.Item(63, 8).Value = $ws.Cells.Item(64, 8).Value()
$ws.Cells.Item(63, 9).Value = $ws.Cells.Item(64, 9).Value()
$ws.Cells.Item(63, 10).Value = $ws.Cells.Item(64, 10).Value()
$ws.Cells.Item(63, 11).Value = $ws.Cells.Item(64, 11).Value()
$ws.Cells.Item(63, 12).Value = $ws.Cells.Item(64, 12).Value()
$ws.Cells.Item(63, 13).Value = 150
$ws.Cells.Item(63, 14).Value = 27000
$ws.Cells.Item(63, 15).Value = 28000
$ws.Cells.Item(63, 16).Value = 27500
$ws.Cells.Item(63, 17).Value = $ws.Cells.Item(64, 17).Value()
$ws.Cells.Item(63, 18).Value = "Perú"
$ws.Cells.Item(63, 19).Value = 1375
$ws.Cells.Item(63, 20).Value = $ws.Cells.Item(64, 20).Value()
